$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the testcase data: E3 should be "pig" instead of the previous value.
$ws.Range("E3").Value = "pig"
